$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Self-evaluation: mark two more criteria as done (column D) in the
# "TESTS ET VALIDATION" section (rows 37 and 38).
$ws.Range("D37").Value = 1
$ws.Range("D38").Value = 1

# The sums in D44 (section subtotal) and D58 (grand total) are formulas
# and will recalculate automatically to reflect the new values.

# Move the view / selection to the newly-updated cell, like the author
# did after entering the new evaluation marks.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D37").Select()
